# Weekly data refresh for the "Ajo" (garlic) sheet.
# A new price observation (week) is inserted as a new row right before the
# existing row 144, pushing that row and every row after it down by one.
# The row that used to be the very last data row (256) ends up at row 257,
# which is why the sheet's used range grows from A1:R256 to A1:R257.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 144; Excel shifts rows 144:256 down to
# 145:257 automatically (formats, dimension, etc. all move with them).
$ws.Rows.Item(144).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(144, 1).Value  = 8
$ws.Cells.Item(144, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(144, 3).Value  = "Coquimbo"
$ws.Cells.Item(144, 4).Value  = 44673
$ws.Cells.Item(144, 5).Value  = 4
$ws.Cells.Item(144, 6).Value  = 100112003
$ws.Cells.Item(144, 7).Value  = "Ajo"
$ws.Cells.Item(144, 8).Value  = "Chino"
$ws.Cells.Item(144, 9).Value  = "Primera"
$ws.Cells.Item(144, 10).Value = 560
$ws.Cells.Item(144, 11).Value = 18500
$ws.Cells.Item(144, 12).Value = 19000
$ws.Cells.Item(144, 13).Value = 18750
$ws.Cells.Item(144, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(144, 15).Value = "China"
$ws.Cells.Item(144, 16).Value = 1875
$ws.Cells.Item(144, 17).Value = 10
$ws.Cells.Item(144, 18).Value = "Hortaliza"
